$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.382.79'
$ws.Range("E2").Value = '  +2.39%  '

$ws.Range("D3").Value = '2.067.13'
$ws.Range("E3").Value = '  +3.63%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.19'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.61%  '

$ws.Range("E6").Value = '  +2.42%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.35'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.90%  '

$ws.Range("E9").Value = '  +3.22%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '59.03'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.92%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0762'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.85%  '

$ws.Range("E12").Value = '  +2.80%  '

$ws.Range("D13").Value = '2.371.40'
$ws.Range("E13").Value = '  +3.64%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.57'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.25%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.39'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.79%  '

$ws.Range("E16").Value = '  +2.53%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.18'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.94%  '

$ws.Range("D18").Value = '2.073.71'
$ws.Range("E18").Value = '  +3.49%  '

$ws.Range("D19").Value = '37.568.02'
$ws.Range("E19").Value = '  +2.97%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.15'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +16.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.19'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.57%  '

$ws.Range("D22").Value = '0.0₃0814'
$ws.Range("E22").Value = '  +1.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '227.29'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.28%  '

$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("E25").Value = '  +2.17%  '

$ws.Range("E26").Value = '  +1.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '165.32'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.51'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +13.65%  '

$ws.Range("E29").Value = '  +2.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.23'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.24%  '

$ws.Range("E31").Value = '  +0.03%  '

$ws.Range("E32").Value = '  +1.88%  '

$ws.Range("E33").Value = '  +3.34%  '

$ws.Range("E34").Value = '  +3.26%  '

$ws.Range("E35").Value = '  +8.14%  '

$ws.Range("E36").Value = '  +7.31%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.38'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.89%  '

$ws.Range("E38").Value = '  +0.03%  '

$ws.Range("E39").Value = '  +1.40%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.86'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.21%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0981'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.58%  '

$ws.Range("E42").Value = '  -1.27%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.44'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +22.82%  '

$ws.Range("D44").Value = '1.459.11'
$ws.Range("E44").Value = '  +0.31%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '95.70'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.42%  '

$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0211'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.35%  '

$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.16'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.80%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.85'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.17%  '

$ws.Range("E49").Value = '  +4.10%  '

$ws.Range("E50").Value = '  +6.55%  '

$ws.Range("E51").Value = '  +1.98%  '
